# Auto-generated Excel COM-interop edit script
# Applies cryptos list update per commit message:
# 'Updated cryptos list on Mon Oct 23 14:47:58 UTC 2023 with GitHub Actions'
#
# Strategy: most cells are plain text (coin names, URLs, price/volume strings)
# and can be written directly via Range.Value. A handful of 'Price' (column D)
# values are strings that LOOK like plain numbers (e.g. "0.999", "66.00") -
# assigning those via .Value would make Excel auto-convert them to a numeric
# cell (and lose the trailing/insignificant zeros the source format relies on).
# To keep them as literal text (matching the workbook's inlineStr cells) without
# leaving a stray NumberFormat/style behind, we briefly mark the cell as Text,
# assign the literal string, then restore the cell's style to Normal so the
# saved file carries no extra formatting - only the value itself changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.687.93'
$ws.Range("E2").Value = '  +2.61%  '
$ws.Range("D3").Value = '1.678.25'
$ws.Range("E3").Value = '  +2.90%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.20%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '219.25'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.18%  '
$ws.Range("E6").Value = '  +1.93%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.21%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '29.11'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.44%  '
$ws.Range("E9").Value = '  +2.08%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0644'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.94%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0904'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.06%  '
$ws.Range("D12").Value = '1.918.54'
$ws.Range("E12").Value = '  +2.87%  '
$ws.Range("D13").Value = '1.675.49'
$ws.Range("E13").Value = '  +2.89%  '
$ws.Range("B14").Value = 'Polygon'
$ws.Range("C14").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.605'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +7.72%  '
$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '10.10'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +8.42%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.02'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.49%  '
$ws.Range("D17").Value = '30.692.89'
$ws.Range("E17").Value = '  +2.61%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '66.00'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.00%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '243.00'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.84%  '
$ws.Range("E20").Value = '  +2.59%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.00'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.06%  '
$ws.Range("E22").Value = '  +2.47%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.96'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.56%  '
$ws.Range("E24").Value = '  -0.65%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '159.06'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.71%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.81'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.40%  '
$ws.Range("E27").Value = '  +2.39%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.69'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.99%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.22%  '
$ws.Range("E30").Value = '  +0.93%  '
$ws.Range("E31").Value = '  +4.25%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.46'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.71%  '
$ws.Range("B33").Value = 'Maker'
$ws.Range("C33").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D33").Value = '1.514.14'
$ws.Range("E33").Value = '  +6.32%  '
$ws.Range("B34").Value = 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.31'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.36%  '
$ws.Range("B36").Value = 'Aave'
$ws.Range("C36").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '83.33'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +10.41%  '
$ws.Range("B37").Value = 'TrustWalletToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.02'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.42%  '
$ws.Range("E38").Value = '  +8.28%  '
$ws.Range("E39").Value = '  +4.57%  '
$ws.Range("E40").Value = '  -2.94%  '
$ws.Range("E41").Value = '  +0.04%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.02'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.84%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.837'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.15%  '
$ws.Range("E44").Value = '  -0.09%  '
$ws.Range("E45").Value = '  +1.79%  '
$ws.Range("E46").Value = '  -0.09%  '
$ws.Range("E47").Value = '  +4.08%  '
$ws.Range("D48").Value = '1.811.07'
$ws.Range("E48").Value = '  +2.17%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '49.89'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.61%  '
$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '92.74'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.67%  '
$ws.Range("B51").Value = 'BabyDogeCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D51").Value = '0.0₆0116'
$ws.Range("E51").Value = '  +2.98%  '
